# Apply the attendance-report sync edits described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Global swap of "Recorded By" ordering: "System, dnasr281@gmail.com"
#    -> "dnasr281@gmail.com, System" everywhere it appears on the sheet.
# ---------------------------------------------------------------------------
$used = $ws.UsedRange
$used.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")

# ---------------------------------------------------------------------------
# 2) Top summary box (K/L columns, rows 6-10) for group B1A1.
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 231          # Recorded Sessions
$ws.Range("L7").Value = 27           # Missing Sessions

# L9 holds a "NN.N%" value stored as literal text (not a number). Writing a
# "...%" string straight into .Value gets auto-parsed into a percentage
# number + reformats the cell's style, so force text with a leading
# apostrophe and then restore the original (fill/font) style by pasting
# just the formats back in from an untouched donor cell of the same style.
$ws.Range("L9").Value = "'72.6%"     # Coverage %
$ws.Range("L10").Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3) Per-group statistics table (rows 15-20, columns K-S).
#    Columns: K=Year L=Group M=Students N=Total O=Recorded P=Missing
#             Q=Pending R=Coverage% S=AvgAttendance%
#    O/P are real numbers (safe to assign directly); R/S are "NN.N%" text
#    values, so they need the same text-preserving dance as L9 above.
# ---------------------------------------------------------------------------

# Row 15 - B1A1
$ws.Range("O15").Value = 20
$ws.Range("P15").Value = 2
$ws.Range("R15").Value = "'76.9%"
$ws.Range("S15").Value = "'80.8%"

# Row 16 - B1A2
$ws.Range("O16").Value = 21
$ws.Range("P16").Value = 1
$ws.Range("R16").Value = "'80.8%"
$ws.Range("S16").Value = "'81.1%"

# Row 17 - B1B1
$ws.Range("O17").Value = 21
$ws.Range("P17").Value = 1
$ws.Range("R17").Value = "'80.8%"
$ws.Range("S17").Value = "'72.2%"

# Row 18 - B1B2
$ws.Range("O18").Value = 21
$ws.Range("P18").Value = 1
$ws.Range("R18").Value = "'80.8%"
$ws.Range("S18").Value = "'78.1%"

# Row 19 - B1C1
$ws.Range("O19").Value = 21
$ws.Range("P19").Value = 1
$ws.Range("R19").Value = "'80.8%"

# Row 20 - B1C2
$ws.Range("O20").Value = 20
$ws.Range("P20").Value = 2
$ws.Range("R20").Value = "'76.9%"
$ws.Range("S20").Value = "'79.3%"

# Restore the original style (lost by the text-coercion above) on every R/S
# cell we just touched, using untouched style-5 donor cells R21/S21.
$ws.Range("R21").Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null
$ws.Range("R16").PasteSpecial(-4122) | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null
$ws.Range("R19").PasteSpecial(-4122) | Out-Null
$ws.Range("R20").PasteSpecial(-4122) | Out-Null

$ws.Range("S21").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null
$ws.Range("S16").PasteSpecial(-4122) | Out-Null
$ws.Range("S17").PasteSpecial(-4122) | Out-Null
$ws.Range("S18").PasteSpecial(-4122) | Out-Null
$ws.Range("S20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Six "11/01/2026" sessions flip from Not Recorded -> Recorded.
#    Copy the (green "Recorded") formatting from the row above onto A:I of
#    each target row, then write the new Recorded-by / Students / Status
#    values.
# ---------------------------------------------------------------------------

$ws.Range("A22:I22").Copy() | Out-Null
$ws.Range("A23:I23").PasteSpecial(-4122) | Out-Null
$ws.Range("G23").Value = "dnasr281@gmail.com"
$ws.Range("H23").Value = "20/26"
$ws.Range("I23").Value = "Recorded"

$ws.Range("A48:I48").Copy() | Out-Null
$ws.Range("A49:I49").PasteSpecial(-4122) | Out-Null
$ws.Range("G49").Value = "dnasr281@gmail.com"
$ws.Range("H49").Value = "23/27"
$ws.Range("I49").Value = "Recorded"

$ws.Range("A74:I74").Copy() | Out-Null
$ws.Range("A75:I75").PasteSpecial(-4122) | Out-Null
$ws.Range("G75").Value = "dnasr281@gmail.com"
$ws.Range("H75").Value = "21/26"
$ws.Range("I75").Value = "Recorded"

$ws.Range("A100:I100").Copy() | Out-Null
$ws.Range("A101:I101").PasteSpecial(-4122) | Out-Null
$ws.Range("G101").Value = "dnasr281@gmail.com"
$ws.Range("H101").Value = "18/27"
$ws.Range("I101").Value = "Recorded"

$ws.Range("A126:I126").Copy() | Out-Null
$ws.Range("A127:I127").PasteSpecial(-4122) | Out-Null
$ws.Range("G127").Value = "dnasr281@gmail.com"
$ws.Range("H127").Value = "23/30"
$ws.Range("I127").Value = "Recorded"

$ws.Range("A152:I152").Copy() | Out-Null
$ws.Range("A153:I153").PasteSpecial(-4122) | Out-Null
$ws.Range("G153").Value = "dnasr281@gmail.com"
$ws.Range("H153").Value = "16/23"
$ws.Range("I153").Value = "Recorded"

$excel.CutCopyMode = 0
